$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10

# Leading single-quote forces the "quote prefix" text semantics (same as typing
# a leading apostrophe in Excel) so the existing quotePrefix-bearing cell style
# is reused instead of Excel forking a brand-new (unprefixed) style.

# Row 18 (Thurs, Mar 18) - "Do Before Class" column (C18):
# Remove the "Backwards Design Due" note (it was incorrectly placed here before).
$c18 = "'" + '- `Taxonomy of Questions <taxonomy_of_questions.ipynb>`_'
$ws.Range("C18").Value = $c18

# Row 17 (Tues, Mar 16) - "Do Before Class" column (C17):
# Replace "Angrist and Pischke (MM), Chapter 4" with "Morgan and Winship, Chapter 11"
# and add the "Backwards Design Due" note here (this is the correct due date).
$c17 = "'" + '- Angrist and Pischke (MM), Chapter 3, Sections 3.3 - End' + $nl + '- Morgan and Winship, Chapter 11' + $nl + '- `RDD at Coursera <https://medium.com/coursera-engineering/regression-discontinuity-understanding-the-benefit-of-subtitles-on-coursera-dd82bb25a0f1>`_' + $nl + '- **Backwards Design Due**'
$ws.Range("C17").Value = $c17

# The extra line makes row 17 taller (grew from a 3-line note to a 4-line note).
$ws.Rows(17).RowHeight = 85
